# Changed dates to ISO
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the workbook/project name.
$ws.Name = "covid19_cases_switzerland"

# Re-format the Date column (A2:A8) from the default short-date format
# to an explicit ISO 8601 (yyyy-mm-dd) custom number format.
$ws.Range("A2:A8").NumberFormat = "yyyy\-mm\-dd;@"

# The new date strings are a different width than the old ones, so the
# column needs to be re-sized to fit the content again.
$ws.Columns("A:A").ColumnWidth = 9.59

# Move/restore the active selection to F16 (matches the saved view state).
[void]$ws.Range("F16").Select()
